$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1829268292682927
$ws.Range("C2").Value = 0.551829268292683
$ws.Range("J2").Value = 0.01219512195121951
$ws.Range("P2").Value = 0.1219512195121951
$ws.Range("S2").Value = 0.1310975609756098
$ws.Range("B3").Value = 0.005319148936170213
$ws.Range("C3").Value = 0.03191489361702127
$ws.Range("J3").Value = 0.04787234042553191
$ws.Range("P3").Value = 0.7446808510638298
$ws.Range("S3").Value = 0.1702127659574468
$ws.Range("P4").Value = 0.6545454545454545
$ws.Range("S4").Value = 0.3454545454545455
$ws.Range("B6").Value = 0.1033057851239669
$ws.Range("F6").Value = 0.05785123966942149
$ws.Range("J6").Value = 0.2479338842975207
$ws.Range("O6").Value = 0.004132231404958678
$ws.Range("Q6").Value = 0.1900826446280992
$ws.Range("R6").Value = 0.07851239669421488
$ws.Range("S6").Value = 0.3181818181818182
$ws.Range("B7").Value = 0.1024390243902439
$ws.Range("D7").Value = 0.01463414634146342
$ws.Range("F7").Value = 0.06341463414634146
$ws.Range("J7").Value = 0.1560975609756098
$ws.Range("O7").Value = 0.02926829268292683
$ws.Range("Q7").Value = 0.175609756097561
$ws.Range("R7").Value = 0.05853658536585366
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.08713692946058091
$ws.Range("D8").Value = 0.01659751037344398
$ws.Range("E8").Value = 0.002074688796680498
$ws.Range("F8").Value = 0.05601659751037345
$ws.Range("J8").Value = 0.1369294605809129
$ws.Range("O8").Value = 0.01867219917012448
$ws.Range("Q8").Value = 0.2157676348547718
$ws.Range("R8").Value = 0.07883817427385892
$ws.Range("S8").Value = 0.3879668049792531
$ws.Range("B9").Value = 0.09302325581395349
$ws.Range("D9").Value = 0.02906976744186046
$ws.Range("E9").Value = 0.005813953488372093
$ws.Range("F9").Value = 0.04651162790697674
$ws.Range("J9").Value = 0.1046511627906977
$ws.Range("O9").Value = 0.04069767441860465
$ws.Range("Q9").Value = 0.2267441860465116
$ws.Range("R9").Value = 0.1046511627906977
$ws.Range("S9").Value = 0.3488372093023256
$ws.Range("B10").Value = 0.1136203246294989
$ws.Range("D10").Value = 0.02681721947776994
$ws.Range("E10").Value = 0.0007057163020465773
$ws.Range("F10").Value = 0.08045165843330981
$ws.Range("J10").Value = 0.1263232180663373
$ws.Range("O10").Value = 0.01129146083274524
$ws.Range("Q10").Value = 0.2335920959774171
$ws.Range("R10").Value = 0.07339449541284404
$ws.Range("S10").Value = 0.3338038108680311
$ws.Range("G11").Value = 0.1213114754098361
$ws.Range("J11").Value = 0.07213114754098361
$ws.Range("K11").Value = 0.1508196721311476
$ws.Range("L11").Value = 0.639344262295082
$ws.Range("S11").Value = 0.01639344262295082
$ws.Range("G12").Value = 0.7009803921568627
$ws.Range("J12").Value = 0.25
$ws.Range("L12").Value = 0.0392156862745098
$ws.Range("S12").Value = 0.00980392156862745
$ws.Range("G13").Value = 0.7045454545454546
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.04545454545454546
$ws.Range("F15").Value = 0.008097165991902834
$ws.Range("H15").Value = 0.1740890688259109
$ws.Range("I15").Value = 0.06072874493927125
$ws.Range("J15").Value = 0.3562753036437247
$ws.Range("K15").Value = 0.08097165991902834
$ws.Range("M15").Value = 0.01619433198380567
$ws.Range("O15").Value = 0.08097165991902834
$ws.Range("S15").Value = 0.2226720647773279
$ws.Range("F16").Value = 0.009615384615384616
$ws.Range("H16").Value = 0.1923076923076923
$ws.Range("I16").Value = 0.08653846153846154
$ws.Range("J16").Value = 0.3605769230769231
$ws.Range("K16").Value = 0.1442307692307692
$ws.Range("M16").Value = 0.02884615384615385
$ws.Range("O16").Value = 0.05288461538461538
$ws.Range("S16").Value = 0.125
$ws.Range("F17").Value = 0.02355072463768116
$ws.Range("H17").Value = 0.197463768115942
$ws.Range("I17").Value = 0.06521739130434782
$ws.Range("J17").Value = 0.447463768115942
$ws.Range("K17").Value = 0.09239130434782608
$ws.Range("M17").Value = 0.01449275362318841
$ws.Range("O17").Value = 0.05072463768115942
$ws.Range("S17").Value = 0.108695652173913
$ws.Range("F18").Value = 0.03174603174603174
$ws.Range("H18").Value = 0.164021164021164
$ws.Range("I18").Value = 0.08994708994708994
$ws.Range("J18").Value = 0.4285714285714285
$ws.Range("K18").Value = 0.08994708994708994
$ws.Range("M18").Value = 0.01587301587301587
$ws.Range("O18").Value = 0.07407407407407407
$ws.Range("S18").Value = 0.1058201058201058
$ws.Range("F19").Value = 0.01436552274541101
$ws.Range("H19").Value = 0.209098164405427
$ws.Range("I19").Value = 0.06943335993615324
$ws.Range("J19").Value = 0.3854748603351955
$ws.Range("K19").Value = 0.1133280127693535
$ws.Range("M19").Value = 0.01915403032721468
$ws.Range("N19").Value = 0.001596169193934557
$ws.Range("O19").Value = 0.08379888268156424
$ws.Range("S19").Value = 0.1037509976057462

Write-Host "Applied 107 cell updates"
